# Auto-generated PowerShell Excel COM-interop script
# Applies the cryptos list update per the commit diff

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value (values are written as Text to preserve exact formatting)
$updates = @{
    'D2' = '68.371.27'
    'E2' = '  +1.27%  '
    'D3' = '3.917.37'
    'E3' = '  -0.81%  '
    'E4' = '  +0.00%  '
    'D5' = '485.62'
    'E5' = '  +1.26%  '
    'D6' = '146.08'
    'E6' = '  -1.82%  '
    'E7' = '  -0.84%  '
    'D8' = '0.998'
    'E8' = '  -0.10%  '
    'D9' = '0.734'
    'E9' = '  +0.14%  '
    'D10' = '0.167'
    'E10' = '  +0.03%  '
    'D11' = '0.0000345'
    'E11' = '  -2.00%  '
    'D12' = '43.32'
    'E12' = '  -0.32%  '
    'D13' = '10.82'
    'E13' = '  +3.45%  '
    'D14' = '4.538.58'
    'E14' = '  -0.77%  '
    'D15' = '3.925.12'
    'E15' = '  -0.58%  '
    'D16' = '14.32'
    'E16' = '  -4.55%  '
    'E17' = '  -1.23%  '
    'D18' = '20.08'
    'E18' = '  -0.08%  '
    'E19' = '  -1.27%  '
    'D20' = '68.412.72'
    'E20' = '  +1.13%  '
    'D21' = '434.54'
    'E21' = '  +0.02%  '
    'D22' = '3.50'
    'E22' = '  +3.38%  '
    'D23' = '15.09'
    'E23' = '  +3.66%  '
    'D24' = '88.19'
    'E24' = '  +0.66%  '
    'D25' = '11.28'
    'E25' = '  +15.18%  '
    'D26' = '11.20'
    'E26' = '  +10.13%  '
    'E27' = '  -2.05%  '
    'D28' = '38.04'
    'E28' = '  -1.65%  '
    'D29' = '5.70'
    'E29' = '  +0.10%  '
    'D30' = '713.73'
    'E30' = '  -0.79%  '
    'D31' = '13.79'
    'E31' = '  +2.35%  '
    'E32' = '  -2.20%  '
    'E33' = '  +4.42%  '
    'D34' = '6.22'
    'E34' = '  +15.34%  '
    'D35' = '41.50'
    'E35' = '  -1.92%  '
    'D36' = '0.0₃0876'
    'E36' = '  +3.86%  '
    'D37' = '60.86'
    'E37' = '  +4.35%  '
    'E38' = '  -3.89%  '
    'D39' = '1.00'
    'E39' = '  +0.04%  '
    'D40' = '0.393'
    'E40' = '  +15.99%  '
    'D41' = '0.0488'
    'E41' = '  +2.57%  '
    'D42' = '2.93'
    'E42' = '  +16.02%  '
    'E43' = '  +2.00%  '
    'D44' = '2.98'
    'E44' = '  +5.58%  '
    'B45' = 'ApeXProtocol'
    'C45' = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
    'D45' = '3.39'
    'E45' = '  +5.78%  '
    'B46' = 'Stellar'
    'C46' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D46' = '0.142'
    'E46' = '  -1.63%  '
    'B47' = 'FirstDigitalUSD'
    'C47' = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
    'D47' = '1.00'
    'E47' = '  +0.04%  '
    'D48' = '3.43'
    'E48' = '  -1.56%  '
    'D49' = '2.14'
    'E49' = '  -4.21%  '
    'D50' = '145.23'
    'E50' = '  -2.77%  '
    'D51' = '0.0₆0338'
    'E51' = '  +30.34%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
